$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title textbox: "TextBox 3" (shape index 3) ---
$title = $s.Shapes.Item(3)

$null = $title.TextFrame.TextRange.InsertAfter("`rThomas Fishwick")
$p2 = $title.TextFrame.TextRange.Paragraphs(2,1)
$p2.Font.Size = 28

$p1 = $title.TextFrame.TextRange.Paragraphs(1,1)
$p1.Font.Bold = $true

# Shape.Height is in points; the target height is 1015663 EMU (914400 EMU/in, 12700 EMU/pt)
$title.Height = 1015663 / 12700

# --- Merge split run ("... just using " + "100 trees.") in "TextBox 16" (shape index 16) ---
$tb16 = $s.Shapes.Item(16)
$tr16 = $tb16.TextFrame.TextRange
$lastPara = $tr16.Paragraphs($tr16.Paragraphs().Count, 1)
$fullText = $lastPara.Text
$wholeRun = $lastPara.Characters(1, $fullText.Length)
$wholeRun.Text = $fullText
